# Apply edit described by the commit:
#  - Merge/replace the two comment strings in I20/I21 into a single
#    comment "To split into smaller tasks", moved from I21 to H21.
#  - Clear I20 and I21 text content.
#  - Remove the custom row heights on rows 20 and 21 (reset to default).
#  - Update the active selection to H25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old comment text in I20 and I21, move the new combined
# comment into H21.
$ws.Range("I20").ClearContents()
$ws.Range("I21").ClearContents()

$ws.Range("H21").Value = "To split into smaller tasks"

# Reset custom row heights on rows 20 and 21 back to the sheet default
# (auto-fit, removing the explicit ht="75"/ht="120" row height override).
$ws.Range("A20:I21").Rows.AutoFit()

# Update the selected/active cell shown when the workbook is reopened.
$ws.Range("H25").Select()
